# Apply "added phlu page types" edit to the Neos node type definition sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header / data cells
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Nodetype"
$ws.Range("C1").Value = "Properties"
$ws.Range("D1").Value = "Supertypes"
$ws.Range("E1").Value = "Childnodes (autocreated collections)"
$ws.Range("F1").Value = "Constraints"
$ws.Range("G1").Value = "Abstract"
$ws.Range("H1").Value = "Package"
$ws.Range("I1").Value = "Implementiert"
$ws.Range("A2").Value = 100
$ws.Range("B2").Value = "PHLU.Neos.NodeTypes:Content"
$ws.Range("D2").Value = "TYPO3.Neos:Content"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "PHLU.Neos.NodeTypes"
$ws.Range("I2").Value = 1
$ws.Range("A3").Value = 101
$ws.Range("B3").Value = "PHLU.Neos.NodeTypes:ContentCollection"
$ws.Range("D3").Value = "TYPO3.Neos:ContentCollection"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "PHLU.Neos.NodeTypes"
$ws.Range("I3").Value = 1
$ws.Range("A4").Value = 102
$ws.Range("B4").Value = "PHLU.Neos.NodeTypes:Page"
$ws.Range("D4").Value = "TYPO3.Neos.NodeTypes:Page"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "PHLU.Neos.NodeTypes"
$ws.Range("I4").Value = 1
$ws.Range("A5").Value = 200
$ws.Range("B5").Value = "PHLU.Neos.NodeTypes:TeaserMixin"
$ws.Range("C5").Value = "Teaserheadline, Teasertext, Teaserimage"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "PHLU.Neos.NodeTypes"
$ws.Range("I5").Value = 1
$ws.Range("A6").Value = 300
$ws.Range("B6").Value = "PHLU.Neos.NodeTypes:Event"
$ws.Range("C6").Value = "Date, Location"
$ws.Range("D6").Value = "PHLU.Neos.NodeTypes:TeaserMixin, PHLU.Neos.NodeTypes:Content"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "PHLU.Neos.NodeTypes"
$ws.Range("I6").Value = 1
$ws.Range("A7").Value = 301
$ws.Range("B7").Value = "PHLU.Neos.NodeTypes:Events"
$ws.Range("D7").Value = "PHLU.Neos.NodeTypes:Content"
$ws.Range("E7").Value = "main: PHLU.Neos.NodeTypes:EventsCollection"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "PHLU.Neos.NodeTypes"
$ws.Range("I7").Value = 1
$ws.Range("A8").Value = 302
$ws.Range("B8").Value = "PHLU.Neos.NodeTypes:EventsCollection"
$ws.Range("D8").Value = "PHLU.Neos.NodeTypes:ContentCollection"
$ws.Range("F8").Value = "PHLU.Neos.NodeTypes:Event"
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = "PHLU.Neos.NodeTypes"
$ws.Range("I8").Value = 1
$ws.Range("A9").Value = 303
$ws.Range("B9").Value = "PHLU.Neos.NodeTypes:NewsItem"
$ws.Range("C9").Value = "Date"
$ws.Range("D9").Value = "PHLU.Neos.NodeTypes:TeaserMixin, PHLU.Neos.NodeTypes:Content"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = "PHLU.Neos.NodeTypes"
$ws.Range("A10").Value = 304
$ws.Range("B10").Value = "PHLU.Neos.NodeTypes:Publication"
$ws.Range("C10").Value = "Date"
$ws.Range("D10").Value = "PHLU.Neos.NodeTypes:TeaserMixin, PHLU.Neos.NodeTypes:Content"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = "PHLU.Neos.NodeTypes"
$ws.Range("A11").Value = 400
$ws.Range("B11").Value = "PHLU.Corporate:Page"
$ws.Range("D11").Value = "PHLU.Neos.NodeTypes:Page"
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "PHLU.Corporate"
$ws.Range("I11").Value = 1
$ws.Range("A12").Value = 401
$ws.Range("B12").Value = "PHLU.Corporate:Page.Home"
$ws.Range("D12").Value = "PHLU.Neos.Corporate:Page"
$ws.Range("F12").Value = "PHLU.Corporate:Overview"
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = "PHLU.Corporate"
$ws.Range("I12").Value = 1
$ws.Range("A13").Value = 402
$ws.Range("B13").Value = "PHLU.Corporate::Event"
$ws.Range("D13").Value = "PHLU.Neos.NodeTypes:Event"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = "PHLU.Corporate"
$ws.Range("I13").Value = 1
$ws.Range("A14").Value = 403
$ws.Range("B14").Value = "PHLU.Corporate::Events"
$ws.Range("D14").Value = "PHLU.Neos.NodeTypes:Events"
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = "PHLU.Corporate"
$ws.Range("I14").Value = 1
$ws.Range("A15").Value = 404
$ws.Range("B15").Value = "PHLU.Corporate:ContentCollection"
$ws.Range("D15").Value = "PHLU.Neos.NodeTypes:ContentCollection"
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = "PHLU.Corporate"
$ws.Range("I15").Value = 1
$ws.Range("A16").Value = 405
$ws.Range("B16").Value = "PHLU.Corporate:Content"
$ws.Range("D16").Value = "PHLU.Neos.NodeTypes:Content"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "PHLU.Corporate"
$ws.Range("I16").Value = 1
$ws.Range("A17").Value = 406
$ws.Range("B17").Value = "PHLU.Corporate:Page:Overview.Tiles"
$ws.Range("D17").Value = "PHLU.Neos.Corporate:Page"
$ws.Range("F17").Value = "PHLU.Corporate:Page:Overview.Onepage, PHLU.Corporate:Page:View.Detail"
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = "PHLU.Corporate"
$ws.Range("I17").Value = 1
$ws.Range("A18").Value = 407
$ws.Range("B18").Value = "PHLU.Corporate:Page:Overview.Onepage"
$ws.Range("D18").Value = "PHLU.Neos.Corporate:Page"
$ws.Range("F18").Value = "PHLU.Corporate:Page:View.Detail, PHLU.Corporate:Page:View.Dossier"
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = "PHLU.Corporate"
$ws.Range("I18").Value = 1
$ws.Range("A19").Value = 408
$ws.Range("B19").Value = "PHLU.Corporate:Page:View.Detail"
$ws.Range("D19").Value = "PHLU.Neos.Corporate:Page"
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = "PHLU.Corporate"
$ws.Range("I19").Value = 1
$ws.Range("A20").Value = 409
$ws.Range("B20").Value = "PHLU.Corporate:Page:View.Dossier"
$ws.Range("D20").Value = "PHLU.Neos.Corporate:Page"
$ws.Range("F20").Value = "PHLU.Corporate:Page:View.Detail"
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = "PHLU.Corporate"
$ws.Range("I20").Value = 1
$ws.Range("A21").Value = 500
$ws.Range("B21").Value = "PHLU.Corporate:Page.Service"
$ws.Range("D21").Value = "PHLU.Neos.Corporate:Page"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "PHLU.Corporate"
$ws.Range("I21").Value = 2
$ws.Range("B30").Value = "docs.phlu.ch hinzufügen!!"

# Column width adjustments (closest achievable via ColumnWidth rounding)
$ws.Columns.Item(2).ColumnWidth = 34.41666666666947
$ws.Columns.Item(3).ColumnWidth = 21.416666666669926
$ws.Columns.Item(5).ColumnWidth = 33.250000000003865
$ws.Columns.Item(6).ColumnWidth = 66.91666666666833

# Update selection
$ws.Range("B22").Select()
